# a2.docx - "finished asg 2 q 6" edit
#
# 1. The stray "_GoBack" bookmark that used to sit right after the
#    "(Matlab Response)" heading is removed from there...
# 2. ...because the variable "H" used in the Q6 Matlab snippet was renamed
#    to "transfer" (both where it is assigned: "H = (R./..." and where it
#    is consumed: "im = imag(H);"), and Word re-dropped the cursor's
#    "_GoBack" mark right after the last edit, i.e. immediately after the
#    new "imag(transfer" text (just before the closing ");").

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark (it will be re-added at the
#     new cursor position once the text edit below is done). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: "H = (R./(R+1i.*w.*L+(" -> "transfer = (R./(R+1i.*w.*L+(" ---
$rng = $d.Content
$rng.Find.Execute("H = (R./(R+1i.*w.*L+(", $true, $false, $false, $false, $false, $true, 1, $false, "transfer = (R./(R+1i.*w.*L+(", 2) | Out-Null

# --- Step 3: "imag(H);" -> "imag(transfer);" and drop "_GoBack" right
#     after the newly typed "transfer" (this is where Word leaves the
#     cursor / last-edit mark). ---
$rng2 = $d.Content
$rng2.Find.Execute("imag(H);", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# "imag(H);" -> the "H" is the 6th character (index 5) of the match.
$hStart = $rng2.Start + 5
$hEnd = $hStart + 1
$hRange = $d.Range($hStart, $hEnd)
$hRange.Text = "transfer"

# Re-create "_GoBack" right after the word "transfer" we just inserted,
# i.e. before the trailing ");".
$bmPos = $hStart + [string]"transfer".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
